# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the First_Noticeable_Increase_Index (C), the
# First_Noticeable_Increase_Cumulative_Value (E) and the Pulse_Width (G)
# columns on each of the Step3_DataPts_* sheets for the rows whose earlier
# samples (below the configured threshold) are now zeroed out before the
# rise-point search runs.

$wb = $excel.ActiveWorkbook

$updates = @{
    "Step3_DataPts_0.5" = @(
        @{ Row = 2; C = 87; E = 0.03472607290402394; G = 27 },
        @{ Row = 3; C = 87; E = 0.01458152946316227; G = 36 },
        @{ Row = 5; C = 88; E = 0.03816506112867032; G = 37 },
        @{ Row = 6; C = 87; E = 0.02364508842298046; G = 19 }
    )
    "Step3_DataPts_0.7" = @(
        @{ Row = 2; C = 87; E = 0.03472607290402394; G = 56 },
        @{ Row = 3; C = 87; E = 0.01458152946316227; G = 57 },
        @{ Row = 5; C = 88; E = 0.03816506112867032; G = 58 },
        @{ Row = 6; C = 87; E = 0.02364508842298046; G = 56 }
    )
    "Step3_DataPts_0.8" = @(
        @{ Row = 2; C = 87; E = 0.03472607290402394; G = 72 },
        @{ Row = 3; C = 87; E = 0.01458152946316227; G = 72 },
        @{ Row = 5; C = 88; E = 0.03816506112867032; G = 73 },
        @{ Row = 6; C = 87; E = 0.02364508842298046; G = 68 }
    )
    "Step3_DataPts_0.9" = @(
        @{ Row = 2; C = 87; E = 0.03472607290402394; G = 97 },
        @{ Row = 3; C = 87; E = 0.01458152946316227; G = 97 },
        @{ Row = 5; C = 88; E = 0.03816506112867032; G = 98 },
        @{ Row = 6; C = 87; E = 0.02364508842298046; G = 96 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($rowUpdate in $updates[$sheetName]) {
        $row = $rowUpdate.Row
        $ws.Cells.Item($row, 3).Value = $rowUpdate.C   # Column C - First_Noticeable_Increase_Index
        $ws.Cells.Item($row, 5).Value = $rowUpdate.E   # Column E - First_Noticeable_Increase_Cumulative_Value
        $ws.Cells.Item($row, 7).Value = $rowUpdate.G   # Column G - Pulse_Width
    }
}
